$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-missing "positive tests" value for Apr 6 (row 33) ---
$ws.Range("B33").Value = 1846

# --- Append a brand-new row 34 with the Apr 7 data ---
$ws.Range("A34").Value = 43928

$row34 = @{
    "C"  = 1454
    "D"  = 10
    "F"  = 16
    "H"  = 240
    "J"  = 288
    "L"  = 259
    "N"  = 255
    "P"  = 190
    "R"  = 115
    "T"  = 79
    "V"  = 2
    "X"  = 701
    "Y"  = 747
    "Z"  = 6
    "AA" = 289
    "AB" = 109
    "AC" = 31
    "AD" = 42
    "AE" = 120
    "AF" = 4
    "AG" = 8
    "AH" = 85
    "AI" = 31
    "AJ" = 32
    "AK" = 7
    "AL" = 25
    "AM" = 14
    "AN" = 25
    "AO" = 31
    "AP" = 15
    "AQ" = 750
    "AR" = 20
    "AS" = 14
    "AT" = 5
    "AU" = 22
    "AV" = 1
    "AW" = 11
    "AX" = 1
    "AY" = 1
    "AZ" = 6
    "BA" = 3
    "BB" = 14
    "BC" = 2
    "BD" = 9
    "BE" = 14
    "BG" = 34
    "BH" = 3
    "BI" = 75
}

foreach ($col in $row34.Keys) {
    $ws.Range($col + "34").Value = $row34[$col]
}

# --- Update the view / scroll position to match the saved workbook state ---
$ws.Range("BJ34").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 53
$win.ScrollRow = 1
$win.Left = 14600
$win.Top = 0
$win.Width = 14200
$win.Height = 18000
